$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Fill in the remaining values for row 5 (G5:O5)
$ws.Range("G5").Value = 7461
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 223
$ws.Range("J5").Value = 10
$ws.Range("K5").Value = 2063
$ws.Range("L5").Value = 623
$ws.Range("M5").Value = 599
$ws.Range("N5").Value = 3519
$ws.Range("O5").Value = 3295

# Add new row 6 with full values A6:O6
$ws.Range("A6").Value = 180
$ws.Range("B6").Value = 10000
$ws.Range("C6").Value = 27
$ws.Range("D6").Value = 262
$ws.Range("E6").Value = 17598
$ws.Range("F6").Value = 642
$ws.Range("G6").Value = 18571
$ws.Range("H6").Value = 24
$ws.Range("I6").Value = 343
$ws.Range("J6").Value = 18
$ws.Range("K6").Value = 5645
$ws.Range("L6").Value = 903
$ws.Range("M6").Value = 540
$ws.Range("N6").Value = 7473
$ws.Range("O6").Value = 7106

# Add new row 7 with partial values A7:B7
$ws.Range("A7").Value = 180
$ws.Range("B7").Value = 25000

# Update the active selection on the sheet to D9
$ws.Activate()
$ws.Range("D9").Select()
